$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 3191.76
$ws.Range("I15").Value = 3191.76
$ws.Range("K15").Value = 9575.280000000001
$ws.Range("M15").Value = -9406.280000000001
$ws.Range("H38").Value = 712.8570999999999
$ws.Range("I38").Value = 712.8570999999999
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 2138.5713
$ws.Range("L38").Value = 0
$ws.Range("N38").Value = -1766.5713
$ws.Range("M38").ClearContents()
$ws.Range("H39").Value = 2305.5715
$ws.Range("I39").Value = 1128
$ws.Range("K39").Value = 3384
$ws.Range("M39").Value = -3088
$ws.Range("H41").Value = 8333677
$ws.Range("J41").Value = 478.42856
$ws.Range("L41").Value = 478.42856
$ws.Range("N41").Value = -1358.42856
$ws.Range("H53").Value = 3321.15
$ws.Range("I53").Value = 2976.6365
$ws.Range("J53").Value = 3742.2222
$ws.Range("K53").Value = 2976.6365
$ws.Range("L53").Value = 3742.2222
$ws.Range("M53").Value = -2339.6365
$ws.Range("N53").Value = -5016.2222
$ws.Range("H113").Value = 88241176
$ws.Range("I113").Value = 142859180
$ws.Range("K113").Value = 142859180
$ws.Range("M113").Value = -142855926
$ws.Range("H116").Value = 41673284
$ws.Range("I116").Value = 83337896
$ws.Range("J116").Value = 8670.666999999999
$ws.Range("K116").Value = 83337896
$ws.Range("L116").Value = 8670.666999999999
$ws.Range("M116").Value = -83334454
$ws.Range("N116").Value = -15554.667
$ws.Range("H132").Value = 1320.1305
$ws.Range("I132").Value = 1310.9302
$ws.Range("J132").Value = 1452
$ws.Range("K132").Value = 3932.7906
$ws.Range("L132").Value = 4356
$ws.Range("M132").Value = -1402.7906
$ws.Range("N132").Value = -9416
$ws.Range("H135").Value = 910164.8
$ws.Range("I135").Value = 1112042.1
$ws.Range("K135").Value = 10008378.9
$ws.Range("M135").Value = -10005843.9
$ws.Range("H137").Value = 4275.1816
$ws.Range("I137").Value = 6036.6
$ws.Range("J137").Value = 2807.3333
$ws.Range("K137").Value = 18109.8
$ws.Range("L137").Value = 8421.999899999999
$ws.Range("M137").Value = -15559.8
$ws.Range("N137").Value = -13521.9999
$ws.Range("H138").Value = 2429.7666
$ws.Range("I138").Value = 2447.3635
$ws.Range("J138").Value = 2419.5789
$ws.Range("K138").Value = 7342.0905
$ws.Range("L138").Value = 7258.736699999999
$ws.Range("M138").Value = -2202.0905
$ws.Range("N138").Value = -17538.7367

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2159841.5
$ws.Range("I32").Value = 2197722
$ws.Range("K32").Value = 2197722
$ws.Range("M32").Value = -2197435
$ws.Range("H61").Value = 4955.271
$ws.Range("I61").Value = 2041.3235
$ws.Range("K61").Value = 2041.3235
$ws.Range("M61").Value = -1829.3235
$ws.Range("H74").Value = 71019.53999999999
$ws.Range("I74").Value = 136288.75
$ws.Range("J74").Value = 5750.3335
$ws.Range("K74").Value = 136288.75
$ws.Range("L74").Value = 5750.3335
$ws.Range("M74").Value = -135414.75
$ws.Range("N74").Value = -7498.3335
$ws.Range("H77").Value = 71019.53999999999
$ws.Range("I77").Value = 136288.75
$ws.Range("J77").Value = 5750.3335
$ws.Range("K77").Value = 681443.75
$ws.Range("L77").Value = 28751.6675
$ws.Range("M77").Value = -677075.75
$ws.Range("N77").Value = -37487.6675
$ws.Range("H102").Value = 1517.6666
$ws.Range("I102").Value = 1521.2
$ws.Range("K102").Value = 1521.2
$ws.Range("M102").Value = 100.8
$ws.Range("H132").Value = 7038.1943
$ws.Range("I132").Value = 6095.9443
$ws.Range("K132").Value = 18287.8329
$ws.Range("M132").Value = -15757.8329
$ws.Range("H136").Value = 4955.271
$ws.Range("I136").Value = 2041.3235
$ws.Range("K136").Value = 6123.970499999999
$ws.Range("M136").Value = -3573.970499999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H93").Value = 42000
$ws.Range("J93").Value = 42000
$ws.Range("L93").Value = 42000
$ws.Range("N93").Value = -45744
$ws.Range("H108").Value = 59376
$ws.Range("J108").Value = 59376
$ws.Range("L108").Value = 59376
$ws.Range("N108").Value = -67056

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 342.5
$ws.Range("J22").Value = 300
$ws.Range("L22").Value = 300
$ws.Range("N22").Value = -1000
$ws.Range("H58").Value = 20842482
$ws.Range("I58").Value = 83335540
$ws.Range("K58").Value = 83335540
$ws.Range("M58").Value = -83335337
$ws.Range("H122").Value = 2423.1667
$ws.Range("I122").Value = 1924
$ws.Range("J122").Value = 2779.7144
$ws.Range("K122").Value = 5772
$ws.Range("L122").Value = 8339.143199999999
$ws.Range("M122").Value = -3322
$ws.Range("N122").Value = -13239.1432
$ws.Range("H125").Value = 51598
$ws.Range("J125").Value = 51598
$ws.Range("L125").Value = 51598
$ws.Range("N125").Value = -56518
$ws.Range("H132").Value = 11770430
$ws.Range("I132").Value = 2181.65
$ws.Range("J132").Value = 28582214
$ws.Range("K132").Value = 6544.950000000001
$ws.Range("L132").Value = 85746642
$ws.Range("M132").Value = -4014.950000000001
$ws.Range("N132").Value = -85751702
$ws.Range("H136").Value = 20842482
$ws.Range("I136").Value = 83335540
$ws.Range("K136").Value = 250006620
$ws.Range("M136").Value = -250004070

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 1819.4286
$ws.Range("I113").Value = 1175.6666
$ws.Range("J113").Value = 2978.2
$ws.Range("K113").Value = 3526.9998
$ws.Range("L113").Value = 8934.599999999999
$ws.Range("M113").Value = -1356.9998
$ws.Range("N113").Value = -13274.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 4086
$ws.Range("I80").Value = 2894
$ws.Range("K80").Value = 2894
$ws.Range("M80").Value = -1896
$ws.Range("H83").Value = 4086
$ws.Range("I83").Value = 2894
$ws.Range("K83").Value = 14470
$ws.Range("M83").Value = -9478
$ws.Range("H122").Value = 2900711.5
$ws.Range("I122").Value = 4529631
$ws.Range("J122").Value = 4854.8887
$ws.Range("K122").Value = 13588893
$ws.Range("L122").Value = 14564.6661
$ws.Range("M122").Value = -13586443
$ws.Range("N122").Value = -19464.6661
$ws.Range("H132").Value = 4546.4
$ws.Range("I132").Value = 1339.8
$ws.Range("J132").Value = 7753
$ws.Range("K132").Value = 4019.4
$ws.Range("L132").Value = 23259
$ws.Range("M132").Value = -1489.4
$ws.Range("N132").Value = -28319
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4750.3237
$ws.Range("I40").Value = 4458.8335
$ws.Range("K40").Value = 4458.8335
$ws.Range("M40").Value = -4322.8335
$ws.Range("H68").Value = 3535.5833
$ws.Range("I68").Value = 2158.5557
$ws.Range("K68").Value = 2158.5557
$ws.Range("M68").Value = -1409.5557
$ws.Range("H71").Value = 3535.5833
$ws.Range("I71").Value = 2158.5557
$ws.Range("K71").Value = 10792.7785
$ws.Range("M71").Value = -7048.7785

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 3500
$ws.Range("J41").Value = 3500
$ws.Range("L41").Value = 3500
$ws.Range("N41").Value = -4280
$ws.Range("H62").Value = 6999.2856
$ws.Range("I62").Value = 6599.2
$ws.Range("K62").Value = 6599.2
$ws.Range("M62").Value = -5975.2
$ws.Range("H65").Value = 6999.2856
$ws.Range("I65").Value = 6599.2
$ws.Range("K65").Value = 32996
$ws.Range("M65").Value = -29876
$ws.Range("H132").Value = 35724692
$ws.Range("I132").Value = 71443540
$ws.Range("K132").Value = 214330620
$ws.Range("M132").Value = -214328090
$ws.Range("H136").Value = 35755860
$ws.Range("I136").Value = 125001520
$ws.Range("J136").Value = 57596.4
$ws.Range("K136").Value = 375004560
$ws.Range("L136").Value = 172789.2
$ws.Range("M136").Value = -375002010
$ws.Range("N136").Value = -177889.2
